$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataMapasECE")

# Update the youth-congress-candidate percentages for the affected departments
$ws.Range("C5").Value  = 0.09433962264150944    # AREQUIPA
$ws.Range("C6").Value  = 0.08888888888888889    # AYACUCHO
$ws.Range("C7").Value  = 0.06930693069306931    # CAJAMARCA
$ws.Range("C8").Value  = 0.16216216216216217    # CALLAO
$ws.Range("C13").Value = 0.10465116279069768    # JUNIN
$ws.Range("C14").Value = 0.08275862068965517    # LA LIBERTAD
$ws.Range("C16").Value = 0.07221350078492936    # LIMA
$ws.Range("C18").Value = 0.1206896551724138     # LORETO
$ws.Range("C20").Value = 0.1282051282051282     # MOQUEGUA
$ws.Range("C22").Value = 0.09174311926605505    # PIURA
$ws.Range("C23").Value = 0.11235955056179775    # PUNO
$ws.Range("C24").Value = 0.07017543859649122    # SAN MARTIN
$ws.Range("C26").Value = 0.08695652173913043    # TUMBES
$ws.Range("C27").Value = 0.10526315789473684    # UCAYALI

# Restore the last active selection recorded in the sheet (E17)
$ws.Range("E17").Select()
